$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 259
$ws.Range("F7").Value = 13078
$ws.Range("F8").Value = 59
$ws.Range("F9").Value = 119
$ws.Range("F10").Value = 264
$ws.Range("F11").Value = 3145
$ws.Range("F13").Value = 6555
$ws.Range("F16").Value = 3458
$ws.Range("F20").Value = 39
$ws.Range("F22").Value = 123
$ws.Range("F23").Value = 46
$ws.Range("F24").Value = 3634
$ws.Range("F25").Value = 96
$ws.Range("F27").Value = 2965
$ws.Range("F28").Value = 415
$ws.Range("F29").Value = 1897
$ws.Range("F31").Value = 220
$ws.Range("F32").Value = 6722
$ws.Range("F34").Value = 1158
$ws.Range("F35").Value = 1998
$ws.Range("F36").Value = 1295
$ws.Range("F37").Value = 104
$ws.Range("F38").Value = 1048
$ws.Range("F40").Value = 215
$ws.Range("F41").Value = 225
$ws.Range("F42").Value = 1151
$ws.Range("F43").Value = 1144
$ws.Range("F45").Value = 1211
$ws.Range("F46").Value = 1796
$ws.Range("F48").Value = 161
$ws.Range("F49").Value = 1174
# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 102
$ws.Range("F17").Value = 7
# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 448
$ws.Range("F3").Value = 616
$ws.Range("F4").Value = 23
# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 448
$ws.Range("F7").Value = 616
$ws.Range("F8").Value = 259
$ws.Range("F9").Value = 13078
$ws.Range("F10").Value = 59
$ws.Range("F11").Value = 119
$ws.Range("F13").Value = 264
$ws.Range("F14").Value = 3146
$ws.Range("F16").Value = 3458
$ws.Range("F19").Value = 39
$ws.Range("F22").Value = 123
$ws.Range("F23").Value = 46
$ws.Range("F24").Value = 3634
$ws.Range("F27").Value = 2965
$ws.Range("F28").Value = 2966
$ws.Range("F29").Value = 415
$ws.Range("F30").Value = 1897
$ws.Range("F32").Value = 220
$ws.Range("F33").Value = 6722
$ws.Range("F34").Value = 102
$ws.Range("F36").Value = 1158
$ws.Range("F37").Value = 1998
$ws.Range("F39").Value = 1295
$ws.Range("F40").Value = 104
$ws.Range("F41").Value = 1048
$ws.Range("F42").Value = 215
$ws.Range("F43").Value = 225
$ws.Range("F44").Value = 1151
$ws.Range("F45").Value = 1211
$ws.Range("F47").Value = 1796
$ws.Range("F50").Value = 161

$wb.Save()